$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: paragraph "3.1. Сюжет: ..."
#   Remove "забрали на замок в воздухе и " (the old "в воздухе" phrase,
#   including the surrounding proofErr markup which goes away with the
#   run it lived in) so the sentence reads "...принцессу и мы должны её
#   спасти". The removed text is deleted via a plain Find/Replace
#   (wdReplaceOne) which naturally merges the touched runs into one.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "забрали на замок в воздухе и ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 1)

# ---------------------------------------------------------------------
# Relocate the lone "_GoBack" bookmark: it used to sit between
# "Платформы" and the following space in the "Мир игры" paragraph; it
# now belongs right before "мы должны её спасти" in the paragraph we
# just edited. Adding a bookmark with a name that already exists moves
# it (Word bookmark names are unique), so this both creates the new
# location and removes the old one in a single step.
# ---------------------------------------------------------------------
$rTarget = $d.Content.Duplicate
$null = $rTarget.Find.Execute(
    "мы должны её спасти", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $d.Range($rTarget.Start, $rTarget.Start))

# ---------------------------------------------------------------------
# Change 2: paragraph "...5. Мир игры: ..."
#   The runs "Платформы" and " " (previously separated only by the
#   bookmark that just moved away) must become a single run
#   "Платформы ". Editing text in this engine merges every run across
#   the whole contiguous span unless a bookmark stops it, so we drop
#   two temporary bookmarks around exactly the span we want touched,
#   rewrite that span as one piece of text, then remove the temporary
#   bookmarks again (removing a bookmark alone never re-merges runs).
# ---------------------------------------------------------------------
$rLeft = $d.Content.Duplicate
$null = $rLeft.Find.Execute(
    "Платформы", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$leftPos = $rLeft.Start
$d.Bookmarks.Add("zzTmpLeft", $d.Range($leftPos, $leftPos))

$rRight = $d.Content.Duplicate
$null = $rRight.Find.Execute(
    "из земли", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$rightPos = $rRight.Start
$d.Bookmarks.Add("zzTmpRight", $d.Range($rightPos, $rightPos))

# The span currently reads "Платформы" + " " (i.e. already the right
# characters, just split across two runs) so writing that same text
# back verbatim would be a text-identical no-op that the engine skips
# re-serializing. Route through a distinct placeholder first so the
# write is recognised as a real change, then set the final text.
$mergeRange = $d.Range($leftPos, $rightPos)
$mergeRange.Text = "zzPlaceholderzz"
$placeholderRange = $d.Range($leftPos, $leftPos + 15)
$placeholderRange.Text = "Платформы "

$d.Bookmarks.Item("zzTmpLeft").Delete()
$d.Bookmarks.Item("zzTmpRight").Delete()
